$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (player name, position(s), team) replacing rows 2-18.
$names     = @("Dennis Schröder","KJ Simpson","Miles Bridges","Scottie Barnes","Brook Lopez","Nikola Vucevic","Shaedon Sharpe","De'Aaron Fox","Ja Morant","Mikal Bridges","Isaiah Collier","Tyler Herro","Josh Giddey","DeMar DeRozan","Evan Mobley","Luka Doncic","P.J. Washington")
$positions = @("PG,SG","PG","SF,PF","PG,SG,SF,PF","C","PF,C","SG,SF","PG","PG","SG,SF,PF","PG","PG,SG","PG,SG,SF","SF,PF","PF,C","PG,SG","SF,PF")
$teams     = @("Detroit Pistons","Charlotte Hornets","Charlotte Hornets","Toronto Raptors","Milwaukee Bucks","Chicago Bulls","Portland Trail Blazers","San Antonio Spurs","Memphis Grizzlies","New York Knicks","Utah Jazz","Miami Heat","Chicago Bulls","Sacramento Kings","Cleveland Cavaliers","Los Angeles Lakers","Dallas Mavericks")

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $positions[$i]
    $ws.Cells.Item($row, 3).Value = $teams[$i]
}
